$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the folder preview/md file entry (doc_id) for the "10-tourisme" folder row
$ws.Range("A2").Value = "tourisme_exemple"
